# Atualização de bases das ligas, do dia: 10-06-2024 às 21:53
# Swap the full data (columns B through AD) between each of the following
# row pairs. Column A (the sequential id) is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    @(61, 62),
    @(156, 157),
    @(228, 229),
    @(305, 306)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    # Columns B (2) through AD (30)
    for ($col = 2; $col -le 30; $col++) {
        $cell1 = $ws.Cells.Item($r1, $col)
        $cell2 = $ws.Cells.Item($r2, $col)

        $v1 = $cell1.Value2
        $v2 = $cell2.Value2

        $cell1.Value2 = $v2
        $cell2.Value2 = $v1
    }
}
